$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.69
$ws.Range("G2").Value = 1.84
$ws.Range("H2").Value = 6.2
$ws.Range("J2").Value = 3.15
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 2.44
$ws.Range("P2").Value = 1.48
$ws.Range("Q2").Value = 2.68
$ws.Range("S2").Value = 5.5
$ws.Range("T2").Value = 2.4
$ws.Range("U2").Value = 1.58
$ws.Range("W2").Value = 2.18
$ws.Range("X2").Value = 10
$ws.Range("AB2").Value = 6.6
$ws.Range("AD2").Value = 34
$ws.Range("AF2").Value = 10.5
$ws.Range("AH2").Value = 1000
$ws.Range("F3").Value = 3.1
$ws.Range("G3").Value = 3.55
$ws.Range("H3").Value = 2.24
$ws.Range("I3").Value = 2.46
$ws.Range("L3").Value = 1.27
$ws.Range("O3").Value = 1.26
$ws.Range("P3").Value = 2.08
$ws.Range("S3").Value = 2.86
$ws.Range("U3").Value = 2.3
$ws.Range("V3").Value = 1.68
$ws.Range("W3").Value = 1.39
$ws.Range("Y3").Value = 14.5
$ws.Range("AA3").Value = 38
$ws.Range("AB3").Value = 17.5
$ws.Range("AE3").Value = 28
$ws.Range("AF3").Value = 29
$ws.Range("AG3").Value = 16.5
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 42
$ws.Range("J4").Value = 2.76
$ws.Range("M4").Value = 1.08
$ws.Range("S4").Value = 3.7
$ws.Range("V4").Value = 1.46
$ws.Range("Y4").Value = 12
$ws.Range("AB4").Value = 12.5
$ws.Range("F5").Value = 1.59
$ws.Range("G5").Value = 1.71
$ws.Range("H5").Value = 7.4
$ws.Range("I5").Value = 9.800000000000001
$ws.Range("J5").Value = 3.4
$ws.Range("K5").Value = 3.95
$ws.Range("L5").Value = 1.48
$ws.Range("M5").Value = 1.12
$ws.Range("N5").Value = 2.48
$ws.Range("O5").Value = 1.54
$ws.Range("P5").Value = 1.49
$ws.Range("Q5").Value = 2.66
$ws.Range("R5").Value = 1.17
$ws.Range("U5").Value = 1.53
$ws.Range("V5").Value = 1.11
$ws.Range("W5").Value = 2.4
$ws.Range("Y5").Value = 22
$ws.Range("AC5").Value = 10.5
$ws.Range("AD5").Value = 1000
$ws.Range("AF5").Value = 8.199999999999999
$ws.Range("AH5").Value = 42
$ws.Range("AJ5").Value = 20
$ws.Range("H6").Value = 1.84
$ws.Range("K6").Value = 4
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 3.35
$ws.Range("O6").Value = 1.37
$ws.Range("Q6").Value = 2.08
$ws.Range("T6").Value = 1.95
$ws.Range("W6").Value = 1.24
$ws.Range("Y6").Value = 8
$ws.Range("AL6").Value = 80
$ws.Range("F7").Value = 1.22
$ws.Range("G7").Value = 1.36
$ws.Range("H7").Value = 17.5
$ws.Range("I7").Value = 28
$ws.Range("J7").Value = 4.8
$ws.Range("K7").Value = 7.8
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 3.25
$ws.Range("O7").Value = 1.34
$ws.Range("P7").Value = 1.8
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.3
$ws.Range("S7").Value = 3.7
$ws.Range("T7").Value = 2.84
$ws.Range("U7").Value = 1.43
$ws.Range("V7").Value = 1.04
$ws.Range("W7").Value = 3.75
$ws.Range("X7").Value = 17
$ws.Range("Y7").Value = 50
$ws.Range("AB7").Value = 7
$ws.Range("AC7").Value = 17.5
$ws.Range("AF7").Value = 7.2
$ws.Range("AG7").Value = 15
$ws.Range("AJ7").Value = 11
$ws.Range("AK7").Value = 23
$ws.Range("AN7").Value = 8.4
$ws.Range("G8").Value = 1.62
$ws.Range("H8").Value = 5.8
$ws.Range("I8").Value = 8.6
$ws.Range("J8").Value = 3.75
$ws.Range("K8").Value = 4.7
$ws.Range("Q8").Value = 1.92
$ws.Range("F9").Value = 1.52
$ws.Range("G9").Value = 1.54
$ws.Range("I9").Value = 7
$ws.Range("K9").Value = 5
$ws.Range("O9").Value = 1.22
$ws.Range("P9").Value = 2.38
$ws.Range("R9").Value = 1.54
$ws.Range("S9").Value = 2.74
$ws.Range("T9").Value = 1.82
$ws.Range("U9").Value = 2.1
$ws.Range("W9").Value = 2.84
$ws.Range("Y9").Value = 27
$ws.Range("AB9").Value = 10
$ws.Range("AH9").Value = 22
$ws.Range("AI9").Value = 80
$ws.Range("AK9").Value = 14.5
$ws.Range("AN9").Value = 6.8
$ws.Range("F10").Value = 1.72
$ws.Range("G10").Value = 1.95
$ws.Range("H10").Value = 5.1
$ws.Range("I10").Value = 6.4
$ws.Range("J10").Value = 3.45
$ws.Range("K10").Value = 4.1
$ws.Range("L10").Value = 1.37
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 3.15
$ws.Range("O10").Value = 1.36
$ws.Range("P10").Value = 1.75
$ws.Range("Q10").Value = 2.06
$ws.Range("R10").Value = 1.28
$ws.Range("S10").Value = 3.4
$ws.Range("T10").Value = 1.94
$ws.Range("U10").Value = 1.85
$ws.Range("V10").Value = 1.18
$ws.Range("W10").Value = 2.04
$ws.Range("X10").Value = 14.5
$ws.Range("Y10").Value = 17
$ws.Range("AB10").Value = 7.8
$ws.Range("AC10").Value = 9.4
$ws.Range("AD10").Value = 23
$ws.Range("AF10").Value = 11
$ws.Range("AG10").Value = 11
$ws.Range("AH10").Value = 24
$ws.Range("AJ10").Value = 21
$ws.Range("AK10").Value = 22
$ws.Range("AN10").Value = 17
$ws.Range("F11").Value = 2.48
$ws.Range("I11").Value = 3.65
$ws.Range("M11").Value = 1.12
$ws.Range("N11").Value = 2.76
$ws.Range("O11").Value = 1.53
$ws.Range("P11").Value = 1.57
$ws.Range("T11").Value = 2.14
$ws.Range("AE11").Value = 55
$ws.Range("AH11").Value = 24
$ws.Range("F12").Value = 1.99
$ws.Range("G12").Value = 2.18
$ws.Range("H12").Value = 4.1
$ws.Range("I12").Value = 5.2
$ws.Range("J12").Value = 3.2
$ws.Range("K12").Value = 3.7
$ws.Range("Q12").Value = 2.08
$ws.Range("F13").Value = 2.24
$ws.Range("H13").Value = 3.9
$ws.Range("I13").Value = 4.2
$ws.Range("J13").Value = 3.1
$ws.Range("K13").Value = 3.2
$ws.Range("L13").Value = 1.57
$ws.Range("M13").Value = 1.11
$ws.Range("N13").Value = 2.68
$ws.Range("O13").Value = 1.55
$ws.Range("P13").Value = 1.55
$ws.Range("Q13").Value = 2.64
$ws.Range("S13").Value = 5.5
$ws.Range("T13").Value = 2.14
$ws.Range("U13").Value = 1.77
$ws.Range("V13").Value = 1.31
$ws.Range("W13").Value = 1.74
$ws.Range("AB13").Value = 7.2
$ws.Range("AE13").Value = 70
$ws.Range("AN13").Value = 32
$ws.Range("H14").Value = 5.2
$ws.Range("I14").Value = 6.8
$ws.Range("L14").Value = 1.41
$ws.Range("R14").Value = 1.33
$ws.Range("S14").Value = 3.1
$ws.Range("U14").Value = 1.94
$ws.Range("AB14").Value = 9.6
$ws.Range("AF14").Value = 12.5
$ws.Range("AG14").Value = 12.5
$ws.Range("AI14").Value = 110
$ws.Range("AM14").Value = 160
$ws.Range("F15").Value = 2.06
$ws.Range("G15").Value = 2.14
$ws.Range("H15").Value = 3.6
$ws.Range("I15").Value = 3.9
$ws.Range("J15").Value = 3.75
$ws.Range("K15").Value = 4.1
$ws.Range("S15").Value = 2.94
$ws.Range("U15").Value = 2.24
$ws.Range("V15").Value = 1.34
$ws.Range("W15").Value = 1.87
$ws.Range("AE15").Value = 980
